# Scale the "value" column (D) from 万元-as-if-base-unit back up by 10000
# (i.e. push/pull scale correction), for data rows 2 through 33.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 33; $row++) {
    $cell = $ws.Cells.Item($row, 4)
    $val = $cell.Value()
    if ($val -ne $null -and $val -ne "") {
        # Use decimal arithmetic (rather than IEEE754 double multiplication)
        # so the rescaled value matches an exact decimal-point shift,
        # avoiding spurious binary floating-point rounding noise. Assigning
        # the decimal's string form lets Excel parse/store it as a number
        # (a direct [decimal]->[double] cast loses precision here).
        $dec = [decimal]$val
        $scaled = $dec * 10000
        $cell.Value = $scaled.ToString()
    }
}
